# Update Excel file via API
#
# 1. Correct three placeholder "9999999" figures on the CollectionUseDelivery
#    sheet (row for FY2019) with the real reported numbers.
# 2. Change which sheet/tab is active when the workbook is opened: move the
#    active tab from "SocialMedia" (last sheet) to "Rankings" (first sheet).

$wb = $excel.ActiveWorkbook

# --- 1. Fix the placeholder values on CollectionUseDelivery ---------------
$ws = $wb.Worksheets.Item("CollectionUseDelivery")
$ws.Range("C2").Value = 4336241   # Article downloads
$ws.Range("D2").Value = 2243861   # E-book downloads
$ws.Range("E2").Value = 0         # Research database searches

# --- 2. Make "Rankings" the active/selected sheet -------------------------
$wb.Worksheets.Item("Rankings").Activate()
